$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna5"
$ws.Cells.Item(2, 3).Value = "Epha2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.230855
$ws.Cells.Item(2, 8).Value = 0.692565
$ws.Cells.Item(2, 9).Value = 0.06377305075821572
$ws.Cells.Item(2, 10).Value = 0.06377305075821572
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 11.08503166666667
$ws.Cells.Item(2, 14).Value = 33.255095
$ws.Cells.Item(2, 15).Value = 0.4259149152633459
$ws.Cells.Item(2, 16).Value = 0.4259149152633459
$ws.Cells.Item(2, 17).Value = 2.559034985408333
$ws.Cells.Item(2, 18).Value = 23.031314868675
$ws.Cells.Item(2, 19).Value = 0.0271618935097705
$ws.Cells.Item(2, 20).Value = 0.0271618935097705

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna5"
$ws.Cells.Item(3, 3).Value = "Epha2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.230855
$ws.Cells.Item(3, 8).Value = 0.692565
$ws.Cells.Item(3, 9).Value = 0.06377305075821572
$ws.Cells.Item(3, 10).Value = 0.06377305075821572
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.312365
$ws.Cells.Item(3, 14).Value = 0.937095
$ws.Cells.Item(3, 15).Value = 0.01200185227312402
$ws.Cells.Item(3, 16).Value = 0.01200185227312402
$ws.Cells.Item(3, 17).Value = 0.072111022075
$ws.Cells.Item(3, 18).Value = 0.648999198675
$ws.Cells.Item(3, 19).Value = 0.0007653947342065446
$ws.Cells.Item(3, 20).Value = 0.0007653947342065446

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna5"
$ws.Cells.Item(4, 3).Value = "Epha2"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.230855
$ws.Cells.Item(4, 8).Value = 0.692565
$ws.Cells.Item(4, 9).Value = 0.06377305075821572
$ws.Cells.Item(4, 10).Value = 0.06377305075821572
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 14.62900266666666
$ws.Cells.Item(4, 14).Value = 43.88700799999999
$ws.Cells.Item(4, 15).Value = 0.5620832324635302
$ws.Cells.Item(4, 16).Value = 0.5620832324635302
$ws.Cells.Item(4, 17).Value = 3.377178410613333
$ws.Cells.Item(4, 18).Value = 30.39460569552
$ws.Cells.Item(4, 19).Value = 0.03584576251423868
$ws.Cells.Item(4, 20).Value = 0.03584576251423868

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efna5"
$ws.Cells.Item(5, 3).Value = "Epha2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.900731333333333
$ws.Cells.Item(5, 8).Value = 8.702194
$ws.Cells.Item(5, 9).Value = 0.8013189515350044
$ws.Cells.Item(5, 10).Value = 0.8013189515350045
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 11.08503166666667
$ws.Cells.Item(5, 14).Value = 33.255095
$ws.Cells.Item(5, 15).Value = 0.4259149152633459
$ws.Cells.Item(5, 16).Value = 0.4259149152633459
$ws.Cells.Item(5, 17).Value = 32.15469868649222
$ws.Cells.Item(5, 18).Value = 289.39228817843
$ws.Cells.Item(5, 19).Value = 0.3412936933419445
$ws.Cells.Item(5, 20).Value = 0.3412936933419446

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna5"
$ws.Cells.Item(6, 3).Value = "Epha2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.900731333333333
$ws.Cells.Item(6, 8).Value = 8.702194
$ws.Cells.Item(6, 9).Value = 0.8013189515350044
$ws.Cells.Item(6, 10).Value = 0.8013189515350045
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.312365
$ws.Cells.Item(6, 14).Value = 0.937095
$ws.Cells.Item(6, 15).Value = 0.01200185227312402
$ws.Cells.Item(6, 16).Value = 0.01200185227312402
$ws.Cells.Item(6, 17).Value = 0.9060869429366667
$ws.Cells.Item(6, 18).Value = 8.154782486430001
$ws.Cells.Item(6, 19).Value = 0.009617311679977746
$ws.Cells.Item(6, 20).Value = 0.009617311679977746

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna5"
$ws.Cells.Item(7, 3).Value = "Epha2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.900731333333333
$ws.Cells.Item(7, 8).Value = 8.702194
$ws.Cells.Item(7, 9).Value = 0.8013189515350044
$ws.Cells.Item(7, 10).Value = 0.8013189515350045
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 14.62900266666666
$ws.Cells.Item(7, 14).Value = 43.88700799999999
$ws.Cells.Item(7, 15).Value = 0.5620832324635302
$ws.Cells.Item(7, 16).Value = 0.5620832324635302
$ws.Cells.Item(7, 17).Value = 42.43480641061689
$ws.Cells.Item(7, 18).Value = 381.913257695552
$ws.Cells.Item(7, 19).Value = 0.4504079465130822
$ws.Cells.Item(7, 20).Value = 0.4504079465130822

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Efna5"
$ws.Cells.Item(8, 3).Value = "Epha2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.4883596666666667
$ws.Cells.Item(8, 8).Value = 1.465079
$ws.Cells.Item(8, 9).Value = 0.1349079977067798
$ws.Cells.Item(8, 10).Value = 0.1349079977067798
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 11.08503166666667
$ws.Cells.Item(8, 14).Value = 33.255095
$ws.Cells.Item(8, 15).Value = 0.4259149152633459
$ws.Cells.Item(8, 16).Value = 0.4259149152633459
$ws.Cells.Item(8, 17).Value = 5.413482369722778
$ws.Cells.Item(8, 18).Value = 48.72134132750499
$ws.Cells.Item(8, 19).Value = 0.05745932841163078
$ws.Cells.Item(8, 20).Value = 0.05745932841163078

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Efna5"
$ws.Cells.Item(9, 3).Value = "Epha2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.4883596666666667
$ws.Cells.Item(9, 8).Value = 1.465079
$ws.Cells.Item(9, 9).Value = 0.1349079977067798
$ws.Cells.Item(9, 10).Value = 0.1349079977067798
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.312365
$ws.Cells.Item(9, 14).Value = 0.937095
$ws.Cells.Item(9, 15).Value = 0.01200185227312402
$ws.Cells.Item(9, 16).Value = 0.01200185227312402
$ws.Cells.Item(9, 17).Value = 0.1525464672783333
$ws.Cells.Item(9, 18).Value = 1.372918205505
$ws.Cells.Item(9, 19).Value = 0.001619145858939725
$ws.Cells.Item(9, 20).Value = 0.001619145858939725

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efna5"
$ws.Cells.Item(10, 3).Value = "Epha2"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.4883596666666667
$ws.Cells.Item(10, 8).Value = 1.465079
$ws.Cells.Item(10, 9).Value = 0.1349079977067798
$ws.Cells.Item(10, 10).Value = 0.1349079977067798
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 14.62900266666666
$ws.Cells.Item(10, 14).Value = 43.88700799999999
$ws.Cells.Item(10, 15).Value = 0.5620832324635302
$ws.Cells.Item(10, 16).Value = 0.5620832324635302
$ws.Cells.Item(10, 17).Value = 7.14421486595911
$ws.Cells.Item(10, 18).Value = 64.297933793632
$ws.Cells.Item(10, 19).Value = 0.07582952343620931
$ws.Cells.Item(10, 20).Value = 0.07582952343620931

